$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")

# Insert a new header row above the existing data and label the columns.
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Price"
$ws.Range("D1").Value = "Stock"
